$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Seed formatting for the new row 6 by copying row 5's cell formats ---
$ws.Range("A5:H5").Copy()
$ws.Range("A6:H6").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Rows.Item(6).RowHeight = 15

# --- New header cell for the "pass" column ---
$ws.Cells.Item(1, 8).Value = "pass"

# --- Row 2: ahostess-test ---
$ws.Cells.Item(2, 1).Value = "ahostess-test"
$ws.Cells.Item(2, 2).Value = "abc"
$ws.Cells.Item(2, 3).Value = "ahostess-test@test.com"
$ws.Cells.Item(2, 4).Value = "973 BRAHMS CT"
$ws.Cells.Item(2, 5).Value = "TROY"
$ws.Cells.Item(2, 6).Value = "Michigan"
$ws.Cells.Item(2, 7).Value = 48085
$ws.Cells.Item(2, 8).Value = "blackdress19"

# --- Row 3: bcohost-test ---
$ws.Cells.Item(3, 1).Value = "bcohost-test"
$ws.Cells.Item(3, 2).Value = "abc"
$ws.Cells.Item(3, 3).Value = "bcohost-test@test.com"
$ws.Cells.Item(3, 4).Value = "974 BRAHMS CT"
$ws.Cells.Item(3, 5).Value = "TROY"
$ws.Cells.Item(3, 6).Value = "Michigan"
$ws.Cells.Item(3, 7).Value = 48085
$ws.Cells.Item(3, 8).Value = "blackdress19"

# --- Row 4: guest1-test ---
$ws.Cells.Item(4, 1).Value = "guest1-test"
$ws.Cells.Item(4, 2).Value = "abc"
$ws.Cells.Item(4, 3).Value = "guest1-test@test.com"
$ws.Cells.Item(4, 4).Value = "975 BRAHMS CT"
$ws.Cells.Item(4, 5).Value = "TROY"
$ws.Cells.Item(4, 6).Value = "Michigan"
$ws.Cells.Item(4, 7).Value = 48085
$ws.Cells.Item(4, 8).Value = "blackdress19"

# --- Row 5: guest2-test ---
$ws.Cells.Item(5, 1).Value = "guest2-test"
$ws.Cells.Item(5, 2).Value = "abc"
$ws.Cells.Item(5, 3).Value = "guest2-test@test.com"
$ws.Cells.Item(5, 4).Value = "976 BRAHMS CT"
$ws.Cells.Item(5, 5).Value = "TROY"
$ws.Cells.Item(5, 6).Value = "Michigan"
$ws.Cells.Item(5, 7).Value = 48085
$ws.Cells.Item(5, 8).Value = "blackdress19"

# --- Row 6 (new): guest3-test ---
$ws.Cells.Item(6, 1).Value = "guest3-test"
$ws.Cells.Item(6, 2).Value = "abc"
$ws.Cells.Item(6, 3).Value = "guest3-test@test.com"
$ws.Cells.Item(6, 4).Value = "976 BRAHMS CT"
$ws.Cells.Item(6, 5).Value = "TROY"
$ws.Cells.Item(6, 6).Value = "Michigan"
$ws.Cells.Item(6, 7).Value = 48085
$ws.Cells.Item(6, 8).Value = "blackdress19"

# --- Refresh mailto hyperlinks on column C for the new e-mail addresses ---
$ws.Range("C2:C6").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:ahostess-test@test.com")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:bcohost-test@test.com")
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:guest1-test@test.com")
$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:guest2-test@test.com")
$ws.Hyperlinks.Add($ws.Range("C6"), "mailto:guest3-test@test.com")
# Adding a hyperlink rewrites the cell style; restore the sheet's normal
# Hyperlink cell style so it matches the original formatting.
$ws.Range("C2:C6").Style = "Hyperlink"

# --- Column widths for the newly introduced zip (G) and pass (H) columns ---
$ws.Columns("G").ColumnWidth = 5.166666666666667
$ws.Columns("H").ColumnWidth = 11.166666666666666

# --- Select the whole sheet, matching the post-edit UI state ---
$ws.Cells.Select()
